# Cab booking select time problem fixed
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CAB")

# Update the CAB booking date string (keeps same quoted-text formatting)
$ws.Range("C2").Value = '"28/07/2021"'

# Update the CAB booking time from 06:30 PM to 06:30 AM
$ws.Range("D2").Value = 0.27083333333333331
